$wb = $excel.ActiveWorkbook

# --- Workbook-level defined names ---
$wb.Names.Add("cellsRange", "=Cells!`$A`$3:`$E`$12")
$wb.Names.Add("check", "=#REF!")
$wb.Names.Add("columnsRange", "=Cells!`$B:`$E")
$wb.Names.Add("rangeD", "=Cells!`$D`$3:`$D`$12")
$wb.Names.Add("rangeE", "=Cells!`$E`$3:`$E`$12")

# --- Sheet "Cells" ---
$ws = $wb.Worksheets.Item("Cells")

$ws.Range("K1").Value = 1

$ws.Range("G3").Formula = "=50+E3"
$ws.Range("H3").Formula = "=AVERAGE(D3:D12)*10"

$ws.Range("G4").Formula = "=50+E4"
$ws.Range("H4").Formula = "=`$E`$4*100"

$ws.Range("G5").Formula = "=50+E5"
$ws.Range("H5").Formula = '=IF(E5<50,"E5","Hello world")'

$ws.Range("G6").Formula = "=50+E6"
$ws.Range("H6").ClearContents()

$ws.Range("G7").Formula = "=50+E7"
$ws.Range("H7").ClearContents()

$ws.Range("G8").Formula = "=50+E8"
$ws.Range("H8").ClearContents()

$ws.Range("G9").Formula = "=50+E9"
$ws.Range("H9").ClearContents()

$ws.Range("G10").Formula = "=50+E10"
$ws.Range("H10").ClearContents()

$ws.Range("G11").Formula = "=50+E11"
$ws.Range("H11").ClearContents()

$ws.Range("G12").Formula = "=50+E12"
$ws.Range("H12").ClearContents()

$ws.Range("F14").Formula = "=SUM(cellsRange)"
$ws.Range("F15").Formula = "=MergedCells!G11+G12"
$ws.Range("F16").Formula = "=F15"

# update the view/selection to match the edited area
$ws.Activate() | Out-Null
$ws.Range("H6").Select() | Out-Null
$excel.ActiveWindow.TopLeftCell = $ws.Range("D1")

# --- Sheet "MergedCells" ---
$ws2 = $wb.Worksheets.Item("MergedCells")

# column H -> G shift for the value next to the C9:F9 merge
$ws2.Range("G9").Value = 11
$ws2.Range("H9").ClearContents()

# new text kept in the merged block, plus filler cells spanning the
# existing B14:E21 merge area (mirrors the style footprint already used
# by the other merged ranges on this sheet)
$ws2.Range("F14").Value = "This should be kept after removing the column"

foreach ($r in 14..21) {
    foreach ($col in @("B", "C", "D", "E")) {
        $addr = "$col$r"
        if ($addr -ne "B14") {
            $ws2.Range($addr).Borders.Item(1).LineStyle = -4142
        }
    }
}

$ws2.Activate() | Out-Null
$ws2.Range("E24").Select() | Out-Null
$excel.ActiveWindow.TopLeftCell = $ws2.Range("A4")

$ws.Activate() | Out-Null
